$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values (supplemental figures update) for rows 2-21, columns B-K
$data = New-Object 'object[,]' 20,10
# Row 2
$data[0,0] = -18.84244549351818
$data[0,1] = 1.93309164890677
$data[0,2] = -18.84244549351818
$data[0,3] = -18.84244549351818
$data[0,4] = -18.84244549351818
$data[0,5] = -18.84244549351818
$data[0,6] = -18.84244549351818
$data[0,7] = -18.84244549351818
$data[0,8] = -18.84244549351818
$data[0,9] = -18.84244549351818
# Row 3
$data[1,0] = -18.84244549351818
$data[1,1] = -18.84244549351818
$data[1,2] = -18.84244549351818
$data[1,3] = -18.84244549351818
$data[1,4] = -18.84244549351818
$data[1,5] = -18.84244549351818
$data[1,6] = -18.84244549351818
$data[1,7] = 1.339600841064875
$data[1,8] = -18.84244549351818
$data[1,9] = -18.84244549351818
# Row 4
$data[2,0] = -18.84244549351818
$data[2,1] = 1.964730948646056
$data[2,2] = 1.672219446335454
$data[2,3] = -18.84244549351818
$data[2,4] = 3.442582160248125
$data[2,5] = -18.84244549351818
$data[2,6] = -18.84244549351818
$data[2,7] = -18.84244549351818
$data[2,8] = 1.124860989128155
$data[2,9] = -18.84244549351818
# Row 5
$data[3,0] = -18.84244549351818
$data[3,1] = 1.68483773054417
$data[3,2] = -18.84244549351818
$data[3,3] = -18.84244549351818
$data[3,4] = -18.84244549351818
$data[3,5] = 2.825338706464547
$data[3,6] = -18.84244549351818
$data[3,7] = -18.84244549351818
$data[3,8] = -18.84244549351818
$data[3,9] = -18.84244549351818
# Row 6
$data[4,0] = -18.84244549351818
$data[4,1] = -18.84244549351818
$data[4,2] = -18.84244549351818
$data[4,3] = -18.84244549351818
$data[4,4] = -18.84244549351818
$data[4,5] = -18.84244549351818
$data[4,6] = -18.84244549351818
$data[4,7] = -18.84244549351818
$data[4,8] = -18.84244549351818
$data[4,9] = -18.84244549351818
# Row 7
$data[5,0] = 2.457545583551707
$data[5,1] = -18.84244549351818
$data[5,2] = -18.84244549351818
$data[5,3] = -18.84244549351818
$data[5,4] = -18.84244549351818
$data[5,5] = -18.84244549351818
$data[5,6] = -18.84244549351818
$data[5,7] = -18.84244549351818
$data[5,8] = -18.84244549351818
$data[5,9] = -18.84244549351818
# Row 8
$data[6,0] = -18.84244549351818
$data[6,1] = -18.84244549351818
$data[6,2] = -18.84244549351818
$data[6,3] = 1.807123832345811
$data[6,4] = -18.84244549351818
$data[6,5] = -18.84244549351818
$data[6,6] = -18.84244549351818
$data[6,7] = -18.84244549351818
$data[6,8] = -18.84244549351818
$data[6,9] = -18.84244549351818
# Row 9
$data[7,0] = 3.858691935022087
$data[7,1] = -18.84244549351818
$data[7,2] = -18.84244549351818
$data[7,3] = -18.84244549351818
$data[7,4] = -18.84244549351818
$data[7,5] = -18.84244549351818
$data[7,6] = -18.84244549351818
$data[7,7] = -18.84244549351818
$data[7,8] = -18.84244549351818
$data[7,9] = -18.84244549351818
# Row 10
$data[8,0] = -18.84244549351818
$data[8,1] = -18.84244549351818
$data[8,2] = -18.84244549351818
$data[8,3] = -18.84244549351818
$data[8,4] = -18.84244549351818
$data[8,5] = -18.84244549351818
$data[8,6] = -18.84244549351818
$data[8,7] = 1.85509989802863
$data[8,8] = -18.84244549351818
$data[8,9] = 2.317406754262351
# Row 11
$data[9,0] = -18.84244549351818
$data[9,1] = -18.84244549351818
$data[9,2] = -18.84244549351818
$data[9,3] = 2.918794469059594
$data[9,4] = -18.84244549351818
$data[9,5] = 2.949108918837661
$data[9,6] = -18.84244549351818
$data[9,7] = -18.84244549351818
$data[9,8] = -18.84244549351818
$data[9,9] = 2.026784599997147
# Row 12
$data[10,0] = -18.84244549351818
$data[10,1] = -18.84244549351818
$data[10,2] = -18.84244549351818
$data[10,3] = -18.84244549351818
$data[10,4] = -18.84244549351818
$data[10,5] = -18.84244549351818
$data[10,6] = -18.84244549351818
$data[10,7] = -18.84244549351818
$data[10,8] = -18.84244549351818
$data[10,9] = -18.84244549351818
# Row 13
$data[11,0] = -18.84244549351818
$data[11,1] = -18.84244549351818
$data[11,2] = -18.84244549351818
$data[11,3] = 2.540284591052302
$data[11,4] = -18.84244549351818
$data[11,5] = -18.84244549351818
$data[11,6] = -18.84244549351818
$data[11,7] = -18.84244549351818
$data[11,8] = 1.714600250928373
$data[11,9] = 1.694276963630145
# Row 14
$data[12,0] = -18.84244549351818
$data[12,1] = -18.84244549351818
$data[12,2] = 1.532752004018338
$data[12,3] = -18.84244549351818
$data[12,4] = -18.84244549351818
$data[12,5] = -18.84244549351818
$data[12,6] = -18.84244549351818
$data[12,7] = -18.84244549351818
$data[12,8] = -18.84244549351818
$data[12,9] = 1.969441073483652
# Row 15
$data[13,0] = -18.84244549351818
$data[13,1] = -18.84244549351818
$data[13,2] = 1.749578046704174
$data[13,3] = -18.84244549351818
$data[13,4] = -18.84244549351818
$data[13,5] = -18.84244549351818
$data[13,6] = -18.84244549351818
$data[13,7] = -18.84244549351818
$data[13,8] = -18.84244549351818
$data[13,9] = -18.84244549351818
# Row 16
$data[14,0] = -18.84244549351818
$data[14,1] = -18.84244549351818
$data[14,2] = -18.84244549351818
$data[14,3] = -18.84244549351818
$data[14,4] = -18.84244549351818
$data[14,5] = -18.84244549351818
$data[14,6] = -18.84244549351818
$data[14,7] = -18.84244549351818
$data[14,8] = 1.912188420234807
$data[14,9] = -18.84244549351818
# Row 17
$data[15,0] = -18.84244549351818
$data[15,1] = 2.149088116258026
$data[15,2] = 1.87356089007262
$data[15,3] = -18.84244549351818
$data[15,4] = -18.84244549351818
$data[15,5] = -18.84244549351818
$data[15,6] = -18.84244549351818
$data[15,7] = 2.207230922830933
$data[15,8] = 2.520500529173011
$data[15,9] = -18.84244549351818
# Row 18
$data[16,0] = -18.84244549351818
$data[16,1] = -18.84244549351818
$data[16,2] = -18.84244549351818
$data[16,3] = -18.84244549351818
$data[16,4] = -18.84244549351818
$data[16,5] = -18.84244549351818
$data[16,6] = -18.84244549351818
$data[16,7] = 1.885472468916906
$data[16,8] = 2.332134630289915
$data[16,9] = -18.84244549351818
# Row 19
$data[17,0] = -18.84244549351818
$data[17,1] = -18.84244549351818
$data[17,2] = 2.034716546723095
$data[17,3] = -18.84244549351818
$data[17,4] = -18.84244549351818
$data[17,5] = -18.84244549351818
$data[17,6] = -18.84244549351818
$data[17,7] = 1.834731308241205
$data[17,8] = -18.84244549351818
$data[17,9] = -18.84244549351818
# Row 20
$data[18,0] = -18.84244549351818
$data[18,1] = 1.08661204515538
$data[18,2] = 1.482935729534656
$data[18,3] = -18.84244549351818
$data[18,4] = 3.190249020923358
$data[18,5] = -18.84244549351818
$data[18,6] = 4.321925179090274
$data[18,7] = 0.9794545817617004
$data[18,8] = -18.84244549351818
$data[18,9] = 1.921568367635295
# Row 21
$data[19,0] = -18.84244549351818
$data[19,1] = 1.323993453458372
$data[19,2] = -18.84244549351818
$data[19,3] = 1.64216993706584
$data[19,4] = -18.84244549351818
$data[19,5] = 2.375579038923746
$data[19,6] = -18.84244549351818
$data[19,7] = -18.84244549351818
$data[19,8] = -18.84244549351818
$data[19,9] = -18.84244549351818

$ws.Range("B2:K21").Value2 = $data
